# Change "thermochemical water splitting" to "hydrocarbon partial oxidation"
# on the RHPF sheet (column header F1 and row label A6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RHPF")

$ws.Range("F1").Value = "hydrocarbon partial oxidation"
$ws.Range("A6").Value = "hydrocarbon partial oxidation"

$ws.Range("F2").Select()
